# Fix typo in Ch 8: correct the "date" column values (column F) on the
# active sheet. Each date was off and needs to be shifted to the
# corrected serial date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44648
$ws.Range("F3").Value = 44647
$ws.Range("F4").Value = 44646
$ws.Range("F5").Value = 44645
$ws.Range("F6").Value = 44644
$ws.Range("F7").Value = 44643
